$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.166.72"
$ws.Range("E2").Value = "  -3.70%  "

$ws.Range("D3").Value = "1.802.66"
$ws.Range("E3").Value = "  -4.02%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'310.43"
$ws.Range("E5").Value = "  -1.92%  "

$ws.Range("D7").Value = "'0.4207"
$ws.Range("E7").Value = "  -2.58%  "

$ws.Range("D8").Value = "'0.3543"
$ws.Range("E8").Value = "  -4.27%  "

$ws.Range("D9").Value = "'0.07103"
$ws.Range("E9").Value = "  -4.33%  "

$ws.Range("D10").Value = "'0.8439"
$ws.Range("E10").Value = "  -4.80%  "

$ws.Range("D11").Value = "'20.11"
$ws.Range("E11").Value = "  -5.16%  "

$ws.Range("D12").Value = "1.911.67"
$ws.Range("E12").Value = "  -0.36%  "

$ws.Range("D13").Value = "'5.319"
$ws.Range("E13").Value = "  -3.26%  "

$ws.Range("D14").Value = "'0.06916"
$ws.Range("E14").Value = "  -0.84%  "

$ws.Range("D15").Value = "'6.348"
$ws.Range("E15").Value = "  -4.17%  "

$ws.Range("D16").Value = "'1.008"
$ws.Range("E16").Value = "  +0.18%  "

$ws.Range("D17").Value = "'80.78"
$ws.Range("E17").Value = "  -0.57%  "

$ws.Range("D18").Value = "'0.000008779"
$ws.Range("E18").Value = "  -3.98%  "

$ws.Range("D19").Value = "'1.003"
$ws.Range("E19").Value = "  -0.07%  "

$ws.Range("D20").Value = "'15.06"
$ws.Range("E20").Value = "  -3.65%  "

$ws.Range("D21").Value = "27.146.88"
$ws.Range("E21").Value = "  -3.63%  "

$ws.Range("D22").Value = "'5.070"
$ws.Range("E22").Value = "  -0.46%  "

$ws.Range("D23").Value = "'10.85"
$ws.Range("E23").Value = "  -1.00%  "

$ws.Range("D24").Value = "2.028.21"
$ws.Range("E24").Value = "  -4.91%  "

$ws.Range("D25").Value = "'1.959"
$ws.Range("E25").Value = "  -1.11%  "

$ws.Range("D26").Value = "'153.42"
$ws.Range("E26").Value = "  -0.74%  "

$ws.Range("D27").Value = "'18.18"
$ws.Range("E27").Value = "  -3.01%  "

$ws.Range("D28").Value = "'5.034"
$ws.Range("E28").Value = "  -7.38%  "

$ws.Range("D29").Value = "'112.82"
$ws.Range("E29").Value = "  -4.85%  "

$ws.Range("D30").Value = "'1.712"
$ws.Range("E30").Value = "  -10.07%  "

$ws.Range("D31").Value = "'0.08887"
$ws.Range("E31").Value = "  -1.05%  "

$ws.Range("D34").Value = "'4.457"
$ws.Range("E34").Value = "  -5.12%  "

$ws.Range("D35").Value = "'1.098"
$ws.Range("E35").Value = "  -6.32%  "

$ws.Range("D36").Value = "'1.005"
$ws.Range("E36").Value = "  +0.12%  "

$ws.Range("D37").Value = "'1.071"
$ws.Range("E37").Value = "  -5.78%  "

$ws.Range("D38").Value = "'0.05207"
$ws.Range("E38").Value = "  -5.02%  "

$ws.Range("D39").Value = "'0.01898"
$ws.Range("E39").Value = "  -3.55%  "

$ws.Range("D40").Value = "'2.750"
$ws.Range("E40").Value = "  -4.73%  "

$ws.Range("D41").Value = "'0.1635"
$ws.Range("E41").Value = "  -3.76%  "

$ws.Range("D42").Value = "'0.4974"
$ws.Range("E42").Value = "  -3.91%  "

$ws.Range("D43").Value = "'6.282"
$ws.Range("E43").Value = "  -8.80%  "

$ws.Range("D44").Value = "'8.173"
$ws.Range("E44").Value = "  -4.76%  "

$ws.Range("D45").Value = "'10.27"
$ws.Range("E45").Value = "  -3.10%  "

$ws.Range("D46").Value = "'104.90"
$ws.Range("E46").Value = "  -0.78%  "

$ws.Range("D47").Value = "'1.004"
$ws.Range("E47").Value = "  +0.06%  "

$ws.Range("D48").Value = "'0.06397"
$ws.Range("E48").Value = "  -2.99%  "

$ws.Range("D49").Value = "'0.4563"
$ws.Range("E49").Value = "  -4.20%  "

$ws.Range("D50").Value = "'1.593"
$ws.Range("E50").Value = "  -4.14%  "

$ws.Range("D51").Value = "'62.97"
$ws.Range("E51").Value = "  -3.71%  "

# Row 32: was ImmutableX -> now HuobiToken
$ws.Range("B32").Value = "HuobiToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D32").Value = "'2.960"
$ws.Range("E32").Value = "  -0.64%  "

# Row 33: was HuobiToken -> now ImmutableX
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'0.7395"
$ws.Range("E33").Value = "  -7.04%  "
